$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp values in column A (rows 2-27) to the new
# forecast-aligned datetime serials.
$ws.Cells.Item(2, 1).Value = 45431.916666666664
$ws.Cells.Item(3, 1).Value = 45431.958333333336
$ws.Cells.Item(4, 1).Value = 45432
$ws.Cells.Item(5, 1).Value = 45432.041666666664
$ws.Cells.Item(6, 1).Value = 45432.08333321759
$ws.Cells.Item(7, 1).Value = 45432.124999826388
$ws.Cells.Item(8, 1).Value = 45432.166666435187
$ws.Cells.Item(9, 1).Value = 45432.208333043978
$ws.Cells.Item(10, 1).Value = 45432.249999652777
$ws.Cells.Item(11, 1).Value = 45432.291666261575
$ws.Cells.Item(12, 1).Value = 45432.333332870374
$ws.Cells.Item(13, 1).Value = 45432.374999479165
$ws.Cells.Item(14, 1).Value = 45432.416666087964
$ws.Cells.Item(15, 1).Value = 45432.458332696762
$ws.Cells.Item(16, 1).Value = 45432.499999305554
$ws.Cells.Item(17, 1).Value = 45432.541665914352
$ws.Cells.Item(18, 1).Value = 45432.583332523151
$ws.Cells.Item(19, 1).Value = 45432.624999131942
$ws.Cells.Item(20, 1).Value = 45432.66666574074
$ws.Cells.Item(21, 1).Value = 45432.708332349539
$ws.Cells.Item(22, 1).Value = 45432.74999895833
$ws.Cells.Item(23, 1).Value = 45432.791665567129
$ws.Cells.Item(24, 1).Value = 45432.833332175927
$ws.Cells.Item(25, 1).Value = 45432.874998784719
$ws.Cells.Item(26, 1).Value = 45432.916665393517
$ws.Cells.Item(27, 1).Value = 45432.958332002316

# Scroll the view back to the top (no frozen/offset topLeftCell) and
# select A2:A27, matching the new selection block.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:A27").Select()
